$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '243.28'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-0.30%'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '30.06'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '14.46%'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.130'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-0.13%'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05671'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '1.50%'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '0.76%'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8398'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '2.56%'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8618'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '2.97%'

$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1336'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '0.21%'

$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06910'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-1.13%'

$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.02864'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-0.93%'

$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09381'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.04%'

$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.001522'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.01%'

$ws.Range('B14').Value = 'CoinExToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.04161'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-10.84%'

$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0006007'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.60%'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006056'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-2.07%'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.507'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-4.03%'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.41%'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.132'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-2.29%'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.03248'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '4.39%'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.1295'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-0.29%'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.642'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-3.19%'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001211'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-2.88%'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004318'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-4.04%'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0001180'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '22.85%'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '0.26%'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03713'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '2.10%'

$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1059'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-23.25%'

$ws.Range('B42').Value = 'KickToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.003412'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-44.35%'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-12.51%'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.009683'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '9.32%'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005095'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-4.76%'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.07%'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.09994'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-30.60%'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002716'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '16.32%'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.07%'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.07%'
